$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room for four new entry rows (93-96) right before the
#    "total" row (currently row 96), pushing the three summary rows
#    (total / target / percentage) down from 96-98 to 101-103 (i.e.
#    leaving rows 97-100 empty) - this also carries the custom row
#    height ("ht=14.25") that lives on the total row down with it.
# ------------------------------------------------------------------
$ws.Rows("93:97").Insert()
# The insert leaves a stray placeholder cell in column B of the extra
# spacer row (97) - clear it so the row disappears again.
$ws.Cells.Item(97, 2).Clear()

# ------------------------------------------------------------------
# 2. New diary rows 93-96.
# ------------------------------------------------------------------
# Row 93 - has a date in column A, copy number-format/alignment from
# an existing date cell (A15) so the style matches exactly.
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(93, 1))
$ws.Cells.Item(93, 1).Value = 44586

$ws.Cells.Item(93, 2).Value = 1
$ws.Cells.Item(93, 3).Value = "väärä käyttäjä/omistaja error lisätty /api/users/:id, /profile sivun 'refresh' korjattu"
$ws.Cells.Item(93, 4).Value = "client/api"

$ws.Cells.Item(94, 2).Value = 3
$ws.Cells.Item(94, 3).Value = "Authenticate, Forbidden komponentit tehty, App siivottu, automaattisesti takaisin aikaisemmalle sivulle login jälkeen"
$ws.Cells.Item(94, 4).Value = "client"

$ws.Cells.Item(95, 2).Value = 1
$ws.Cells.Item(95, 3).Value = "PrivateRoute testausta, vanha tapa ei toimi react-router v6:ssa"
$ws.Cells.Item(95, 4).Value = "client"

$ws.Cells.Item(96, 2).Value = 3
$ws.Cells.Item(96, 3).Value = "uusi tapa, wrapper toiminnassa,  refresh saatu taas toimimaan private routen kanssa"
$ws.Cells.Item(96, 4).Value = "client"

# ------------------------------------------------------------------
# 3. Fix up the relocated summary rows (101-103).
# ------------------------------------------------------------------
$ws.Cells.Item(101, 2).Formula = "=SUM(B2:B96)"
$ws.Cells.Item(103, 2).Formula = "=B101/B102*100"

# ------------------------------------------------------------------
# 4. Match the final view/selection state (scrolled down, D97 active).
# ------------------------------------------------------------------
$ws.Range("D97").Select() | Out-Null
